$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.401.82'
$ws.Range('E2').Value = '  +0.43%  '
$ws.Range('D3').Value = '2.381.80'
$ws.Range('E3').Value = '  +0.37%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '550.22'
$ws.Range('E5').Value = '  +0.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.45'
$ws.Range('E6').Value = '  -1.21%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.523'
$ws.Range('E8').Value = '  -0.54%  '
$ws.Range('D9').Value = '2.381.35'
$ws.Range('E9').Value = '  +0.49%  '
$ws.Range('E10').Value = '  +2.49%  '
$ws.Range('E11').Value = '  +1.38%  '
$ws.Range('E12').Value = '  +1.05%  '
$ws.Range('E13').Value = '  +1.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.25'
$ws.Range('E14').Value = '  +0.27%  '
$ws.Range('E15').Value = '  +1.87%  '
$ws.Range('D16').Value = '61.330.99'
$ws.Range('E16').Value = '  +0.31%  '
$ws.Range('D17').Value = '2.376.51'
$ws.Range('E17').Value = '  -0.17%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.96'
$ws.Range('E18').Value = '  +2.61%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '321.53'
$ws.Range('E19').Value = '  +1.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.14'
$ws.Range('E20').Value = '  +0.92%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.75'
$ws.Range('E21').Value = '  +1.64%  '
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '64.31'
$ws.Range('E23').Value = '  +0.93%  '
$ws.Range('E24').Value = '  -9.66%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.61'
$ws.Range('E25').Value = '  +5.68%  '
$ws.Range('E26').Value = '  +1.82%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '511.37'
$ws.Range('E27').Value = '  -2.39%  '
$ws.Range('D28').Value = '0.0₃0897'
$ws.Range('E28').Value = '  -2.56%  '
$ws.Range('E29').Value = '  +3.44%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.37'
$ws.Range('E30').Value = '  -3.04%  '
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('E32').Value = '  -1.39%  '
$ws.Range('E33').Value = '  -0.11%  '
$ws.Range('B34').Value = 'Stacks'
$ws.Range('C34').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.92'
$ws.Range('E34').Value = '  +4.11%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.70'
$ws.Range('E35').Value = '  +1.40%  '
$ws.Range('E36').Value = '  +1.23%  '
$ws.Range('E37').Value = '  +1.21%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.54'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '146.68'
$ws.Range('E39').Value = '  +4.66%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '41.23'
$ws.Range('E41').Value = '  +2.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '150.60'
$ws.Range('E42').Value = '  +7.66%  '
$ws.Range('E43').Value = '  +1.76%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.60'
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('E45').Value = '  +1.51%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.48'
$ws.Range('E46').Value = '  -2.30%  '
$ws.Range('E47').Value = '  +1.14%  '
$ws.Range('E48').Value = '  +0.53%  '
$ws.Range('E49').Value = '  -0.09%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '11.43'
$ws.Range('E50').Value = '  +0.54%  '
$ws.Range('E51').Value = '  +0.60%  '
